$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: strip all whitespace from header labels ---
$ws.Range("A1").Value = "StockCode"
$ws.Range("B1").Value = "CompanyName"
$ws.Range("C1").Value = "Sector"
$ws.Range("D1").Value = "Open"
$ws.Range("E1").Value = "Close"
$ws.Range("F1").Value = "Volume"
$ws.Range("G1").Value = "TradeDate"
$ws.Range("H1").Value = "MarketCap"

# --- Data rows: trim the padded Stock Code / Sector text, and strip all
#     whitespace out of the Company Name text ---
$ws.Range("A2").Value = "BHP"
$ws.Range("B2").Value = "BHPGroupLtd"
$ws.Range("C2").Value = "Materials"

$ws.Range("A3").Value = "CBA"
$ws.Range("B3").Value = "CommonwealthBank"
$ws.Range("C3").Value = "Financials"

$ws.Range("A4").Value = "WBC"
$ws.Range("B4").Value = "WestpacBanking"
$ws.Range("C4").Value = "Financials"

$ws.Range("A5").Value = "CSL"
$ws.Range("B5").Value = "CSLLimited"
$ws.Range("C5").Value = "Healthcare"

$ws.Range("A6").Value = "RIO"
$ws.Range("B6").Value = "RioTintoLimited"
$ws.Range("C6").Value = "Materials"

$ws.Range("A7").Value = "ANZ"
$ws.Range("B7").Value = "ANZBankingGroup"
$ws.Range("C7").Value = "Financials"

$ws.Range("A8").Value = "Joshi"
$ws.Range("B8").Value = "ANZBankingGroup"
$ws.Range("C8").Value = "Financials"

# --- Trade Date column: replace the literal " 2024-03-01 " text with a
#     real date value, formatted as a date (numFmtId 14) ---
$ws.Range("G2").NumberFormat = "mm-dd-yy"
for ($r = 2; $r -le 8; $r++) {
  $ws.Cells.Item($r, 7).Value = [DateTime]"2024-03-01"
}
$ws.Range("G2").Copy()
$ws.Range("G2:G8").PasteSpecial(-4122)

# --- Selection / view state ---
[void]$ws.Range("D10").Select()
